$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1266.3658
$ws.Range("I15").Value = 1266.3658
$ws.Range("K15").Value = 3799.0974
$ws.Range("M15").Value = -3630.0974
$ws.Range("H33").Value = 167391.17
$ws.Range("I33").Value = 366
$ws.Range("K33").Value = 366
$ws.Range("M33").Value = -137
$ws.Range("H53").Value = 221.125
$ws.Range("I53").Value = 211
$ws.Range("K53").Value = 211
$ws.Range("M53").Value = 426
$ws.Range("H62").Value = 3128.7856
$ws.Range("J62").Value = 5155.5713
$ws.Range("L62").Value = 5155.5713
$ws.Range("N62").Value = -6403.5713
$ws.Range("H65").Value = 3128.7856
$ws.Range("J65").Value = 5155.5713
$ws.Range("L65").Value = 25777.8565
$ws.Range("N65").Value = -32017.8565
$ws.Range("H112").Value = 6133.2046
$ws.Range("J112").Value = 6240.9536
$ws.Range("L112").Value = 18722.8608
$ws.Range("N112").Value = -20938.8608
$ws.Range("H115").Value = 829.2
$ws.Range("I115").Value = 829.2
$ws.Range("K115").Value = 2487.6
$ws.Range("M115").Value = -920.6000000000004
$ws.Range("H116").Value = 6017.857
$ws.Range("I116").Value = 7535
$ws.Range("J116").Value = 5765
$ws.Range("K116").Value = 7535
$ws.Range("L116").Value = 5765
$ws.Range("M116").Value = -4093
$ws.Range("N116").Value = -12649
$ws.Range("H118").Value = 1502.4117
$ws.Range("I118").Value = 1492.5714
$ws.Range("K118").Value = 4477.7142
$ws.Range("M118").Value = -2820.7142
$ws.Range("H132").Value = 2767.463
$ws.Range("I132").Value = 2662.4807
$ws.Range("J132").Value = 5497
$ws.Range("K132").Value = 7987.4421
$ws.Range("L132").Value = 16491
$ws.Range("M132").Value = -5457.4421
$ws.Range("N132").Value = -21551
$ws.Range("H139").Value = 88941.96000000001
$ws.Range("J139").Value = 88941.96000000001
$ws.Range("L139").Value = 88941.96000000001
$ws.Range("N139").Value = -99221.96000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9004.671
$ws.Range("I32").Value = 5502.3945
$ws.Range("K32").Value = 5502.3945
$ws.Range("M32").Value = -5215.3945
$ws.Range("H103").Value = 108215.5
$ws.Range("J103").Value = 108215.5
$ws.Range("L103").Value = 108215.5
$ws.Range("N103").Value = -110559.5
$ws.Range("H110").Value = 3592.0476
$ws.Range("J110").Value = 3659.6
$ws.Range("L110").Value = 3659.6
$ws.Range("N110").Value = -7749.6
$ws.Range("H132").Value = 2532.4465
$ws.Range("I132").Value = 1698.4255
$ws.Range("K132").Value = 5095.2765
$ws.Range("M132").Value = -2565.2765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3666.3572
$ws.Range("I20").Value = 3741.3333
$ws.Range("J20").Value = 3531.4
$ws.Range("K20").Value = 3741.3333
$ws.Range("L20").Value = 3531.4
$ws.Range("M20").Value = -3494.3333
$ws.Range("N20").Value = -4025.4
$ws.Range("H94").Value = 857.6316
$ws.Range("I94").Value = 829.7778
$ws.Range("K94").Value = 829.7778
$ws.Range("M94").Value = -378.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7774.1
$ws.Range("I31").Value = 6252.7915
$ws.Range("J31").Value = 9178.385
$ws.Range("K31").Value = 6252.7915
$ws.Range("L31").Value = 9178.385
$ws.Range("M31").Value = -5957.7915
$ws.Range("N31").Value = -9768.385
$ws.Range("H34").Value = 7774.1
$ws.Range("I34").Value = 6252.7915
$ws.Range("J34").Value = 9178.385
$ws.Range("K34").Value = 6252.7915
$ws.Range("L34").Value = 9178.385
$ws.Range("M34").Value = -6050.7915
$ws.Range("N34").Value = -9582.385
$ws.Range("H58").Value = 4021.149
$ws.Range("I58").Value = 4286.6665
$ws.Range("K58").Value = 4286.6665
$ws.Range("M58").Value = -4083.6665
$ws.Range("H69").Value = 10082.5
$ws.Range("I69").Value = 10082.5
$ws.Range("K69").Value = 10082.5
$ws.Range("M69").Value = -9333.5
$ws.Range("H72").Value = 10082.5
$ws.Range("I72").Value = 10082.5
$ws.Range("K72").Value = 30247.5
$ws.Range("M72").Value = -26503.5
$ws.Range("H97").Value = 98197
$ws.Range("J97").Value = 98197
$ws.Range("L97").Value = 98197
$ws.Range("N97").Value = -100179
$ws.Range("H134").Value = 1722.0476
$ws.Range("I134").Value = 1598.2778
$ws.Range("J134").Value = 2464.6667
$ws.Range("K134").Value = 4794.8334
$ws.Range("L134").Value = 7394.000100000001
$ws.Range("M134").Value = -2259.8334
$ws.Range("N134").Value = -12464.0001
$ws.Range("H136").Value = 4021.149
$ws.Range("I136").Value = 4286.6665
$ws.Range("K136").Value = 12859.9995
$ws.Range("M136").Value = -10309.9995
$ws.Range("H141").Value = 215875.05
$ws.Range("J141").Value = 225090.33
$ws.Range("L141").Value = 225090.33
$ws.Range("N141").Value = -235450.33

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2619.5625
$ws.Range("J122").Value = 2868.2964
$ws.Range("L122").Value = 25814.6676
$ws.Range("N122").Value = -30714.6676
$ws.Range("H137").Value = 4737.375
$ws.Range("J137").Value = 4804.8887
$ws.Range("L137").Value = 14414.6661
$ws.Range("N137").Value = -24614.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2751.45
$ws.Range("I102").Value = 2702
$ws.Range("J102").Value = 3031.6667
$ws.Range("K102").Value = 2702
$ws.Range("L102").Value = 3031.6667
$ws.Range("M102").Value = -1080
$ws.Range("N102").Value = -6275.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 482.2
$ws.Range("I93").Value = 467
$ws.Range("K93").Value = 467
$ws.Range("M93").Value = 781
$ws.Range("H125").Value = 89399.92
$ws.Range("J125").Value = 89399.92
$ws.Range("L125").Value = 89399.92
$ws.Range("N125").Value = -99239.92
$ws.Range("H138").Value = 89999.09
$ws.Range("J138").Value = 89999.09
$ws.Range("L138").Value = 89999.09
$ws.Range("N138").Value = -100279.09

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 89675
$ws.Range("J93").Value = 89675
$ws.Range("L93").Value = 89675
$ws.Range("N93").Value = -94667
$ws.Range("H126").Value = 2077.5881
$ws.Range("I126").Value = 2077.5881
$ws.Range("K126").Value = 6232.7643
$ws.Range("M126").Value = -3762.7643
$ws.Range("H136").Value = 14145.155
$ws.Range("I136").Value = 19613.033
$ws.Range("J136").Value = 2037.7142
$ws.Range("K136").Value = 58839.099
$ws.Range("L136").Value = 6113.142599999999
$ws.Range("M136").Value = -56289.099
$ws.Range("N136").Value = -11213.1426
